$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 6925
$ws.Range('I3').Value = 7184
$ws.Range('I4').Value = 1649
$ws.Range('I5').Value = 675
$ws.Range('I6').Value = 8457
$ws.Range('I7').Value = 24890
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('I6').Value = 117
$ws.Range('I7').Value = 289
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('I2').Value = 47
$ws.Range('I6').Value = 50
$ws.Range('I7').Value = 140
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I2').Value = 228
$ws.Range('I7').Value = 767
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I2').Value = 123
$ws.Range('I7').Value = 435
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I2').Value = 224
$ws.Range('I7').Value = 941
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('I3').Value = 64
$ws.Range('I6').Value = 98
$ws.Range('I7').Value = 252
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('I2').Value = 80
$ws.Range('I7').Value = 220
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I6').Value = 180
$ws.Range('I7').Value = 777
$ws.Range('I8').Value = 1484
$ws.Range('I10').Value = 181
$ws.Range('I13').Value = 40
$ws.Range('I14').Value = 140
$ws.Range('I15').Value = 289
$ws.Range('I19').Value = 699
$ws.Range('I20').Value = 616
$ws.Range('I26').Value = 35
$ws.Range('I29').Value = 1488
$ws.Range('I31').Value = 252
$ws.Range('I33').Value = 1096
$ws.Range('I36').Value = 338
$ws.Range('I37').Value = 767
$ws.Range('I39').Value = 19
$ws.Range('I42').Value = 937
$ws.Range('I44').Value = 187
$ws.Range('I45').Value = 48
$ws.Range('I46').Value = 56
$ws.Range('I47').Value = 181
$ws.Range('I48').Value = 315
$ws.Range('I51').Value = 292
$ws.Range('I52').Value = 562
$ws.Range('I55').Value = 287
$ws.Range('I56').Value = 27
$ws.Range('I57').Value = 102
$ws.Range('I60').Value = 143
$ws.Range('I63').Value = 76
$ws.Range('I67').Value = 941
$ws.Range('I73').Value = 226
$ws.Range('I76').Value = 355
$ws.Range('I77').Value = 151
$ws.Range('I78').Value = 332
$ws.Range('I80').Value = 79
$ws.Range('I83').Value = 534
$ws.Range('I84').Value = 220
$ws.Range('I85').Value = 1111
$ws.Range('I86').Value = 160
$ws.Range('I87').Value = 64
$ws.Range('I88').Value = 228
$ws.Range('I90').Value = 324
$ws.Range('I94').Value = 253
$ws.Range('I95').Value = 381
$ws.Range('I96').Value = 289
$ws.Range('I97').Value = 223
$ws.Range('I99').Value = 435
$ws.Range('I101').Value = 24890
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I3').Value = 195
$ws.Range('I7').Value = 534
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('I2').Value = 135
$ws.Range('I6').Value = 77
$ws.Range('I7').Value = 381
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I3').Value = 405
$ws.Range('I6').Value = 352
$ws.Range('I7').Value = 1096
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 440
$ws.Range('I6').Value = 410
$ws.Range('I7').Value = 1488
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I6').Value = 225
$ws.Range('I7').Value = 699
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('I2').Value = 62
$ws.Range('I7').Value = 187
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I2').Value = 51
$ws.Range('I7').Value = 315
$ws = $wb.Worksheets.Item('River North')
$ws.Range('I3').Value = 79
$ws.Range('I6').Value = 162
$ws.Range('I7').Value = 355
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I6').Value = 291
$ws.Range('I7').Value = 1111
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('I6').Value = 52
$ws.Range('I7').Value = 180
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I6').Value = 375
$ws.Range('I7').Value = 937
$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('I5').Value = 15
$ws.Range('I6').Value = 40
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('I6').Value = 84
$ws.Range('I7').Value = 181
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I6').Value = 121
$ws.Range('I7').Value = 332
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('I3').Value = 90
$ws.Range('I7').Value = 287
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('I2').Value = 16
$ws.Range('I3').Value = 19
$ws.Range('I7').Value = 56
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I2').Value = 172
$ws.Range('I3').Value = 174
$ws.Range('I6').Value = 215
$ws.Range('I7').Value = 616
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I3').Value = 113
$ws.Range('I7').Value = 338
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I3').Value = 183
$ws.Range('I6').Value = 183
$ws.Range('I7').Value = 562
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I2').Value = 49
$ws.Range('I6').Value = 146
$ws.Range('I7').Value = 253
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('I6').Value = 60
$ws.Range('I7').Value = 181
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('I3').Value = 69
$ws.Range('I7').Value = 289
$ws = $wb.Worksheets.Item('East Village')
$ws.Range('I2').Value = 8
$ws.Range('I7').Value = 35
$ws = $wb.Worksheets.Item('Greektown')
$ws.Range('I4').Value = 2
$ws.Range('I6').Value = 19
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I5').Value = 3
$ws.Range('I7').Value = 226
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('I6').Value = 145
$ws.Range('I7').Value = 223
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I6').Value = 71
$ws.Range('I7').Value = 228
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 437
$ws.Range('I7').Value = 1484
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('I4').Value = 77
$ws.Range('I6').Value = 37
$ws.Range('I7').Value = 160
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('I6').Value = 115
$ws.Range('I7').Value = 324
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I6').Value = 118
$ws.Range('I7').Value = 292
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('I2').Value = 38
$ws.Range('I3').Value = 25
$ws.Range('I7').Value = 102
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('I2').Value = 51
$ws.Range('I7').Value = 143
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('I2').Value = 50
$ws.Range('I3').Value = 51
$ws.Range('I7').Value = 151
$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range('I3').Value = 9
$ws.Range('I7').Value = 48
$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range('I4').Value = 3
$ws.Range('I7').Value = 27
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('I6').Value = 45
$ws.Range('I7').Value = 79
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I6').Value = 210
$ws.Range('I7').Value = 777
$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('I6').Value = 38
$ws.Range('I7').Value = 64
